$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new Kadastro record (record no. "2") needs to be inserted as the new
# row 2, pushing the previously-existing row 2 record ("1" / Merkez /
# Gökhan ELGÜL) down to row 3.
$ws.Rows.Item(2).Insert()

# The new row's cells must be stored as text (not auto-converted to
# numbers/dates) so values like "2025-07-16" stay literal text.
$newRow = $ws.Range("A2:G2")
$newRow.NumberFormat = "@"

$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "2025-07-16"
$ws.Range("C2").Value = "İlçe"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "Cins D."
$ws.Range("G2").Value = "Göktan ELGÜL"
